$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 50-52: bump "aula" number column (B) from 50 to 51 ---
$ws.Range("B50").Value = 51
$ws.Range("B51").Value = 51
$ws.Range("B52").Value = 51

# --- Seed rows 53 & 54 with row-52's formatting (style ids + column widths) ---
$ws.Range("B52:E52").Copy()
$ws.Range("B53:E54").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Shared-string pool order matches the authored workbook: JPQL note (E54)
# first, then the lesson-53 title (D53, reused by D54), then the
# RequestParam note (E53).
$ws.Range("E54").Value = "5:27`ninteressante: foi demonstrado o tanto de código necessário para efetuar uma consulta JPQL e como evitar isso usando um método ja criado anteriormente na aula 21 onde foi implementada uma classe abstrata e genérica AbstractDao"
$ws.Range("D53").Value = "`n53. Buscar funcionário por nome"
$ws.Range("E53").Value = "`n1:54`nanotação @RequestParam  tem mesma usabilidade que o @PathVariable porém com diferenças:`n@RequestParam faz troca de valores e variaveis entre controller e view através da Request/Requisição de forma CODIFICADA. podemos ter vários parametros passados por url mas que não são parte da url em sí.`n@PathVariable faz troca de valores e variaveis entre controller e view através da path URL e NÃO CODIFICADO. Em resumo, o @PathVariable é utilizado quando o valor da variável é passada diretamente na URL, mas não como um parametro que você passa após o sinal de interrogação (?) mas sim quando o valor faz parte da url"
$ws.Range("D54").Value = "`n53. Buscar funcionário por nome"

# --- Row 53: new entry "53. Buscar funcionário por nome" (title row) ---
$ws.Range("B53").Value = 53
$ws.Range("C53").Value = $ws.Range("C52").Value2
$ws.Rows.Item(53).RowHeight = 165

# --- Row 54: new entry continuing aula 53 ---
$ws.Range("B54").Value = 53
$ws.Range("C54").Value = $ws.Range("C52").Value2
$ws.Rows.Item(54).RowHeight = 60

# --- Update the sheet view to match the post-edit scroll/selection state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 50
$ws.Range("E57").Select()
